# feat: add 2022-Q1 data
#
# Before:
#   Sheet 1 "2021-Q4" - per-fund holdings table for 2021-Q4
#   Sheet 2 "总计"      - summary table (one row per quarter, 2021-Q4 only)
#
# After:
#   Sheet 1 "2021-Q4"  - unchanged
#   Sheet 2 "2022-Q1"  - per-fund holdings table for 2022-Q1 (same layout
#                        as "2021-Q4", built from the old "总计" sheet)
#   Sheet 3 "总计"      - brand-new summary sheet listing both quarters
#                        (2022-Q1 first, then 2021-Q4)

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet into the "2022-Q1" holdings
# sheet. Its B1/C1/D1 header cells and A2 cell already carry the
# existing "header" style used throughout this workbook, so simply
# overwriting their text keeps that formatting intact.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

# E1:H1 are new header cells - give them the same look as B1:D1 by
# copying B1's cell format onto them.
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial($xlPasteFormats)

# A2 already carries the header style from the old sheet - leave its
# format alone, only refresh its value.
$q1.Range("A2").Value = 0

$q1.Range("B2").Value = "004250"
$q1.Range("C2").Value = "银河量化优选混合"

# D2:G2 are numeric-looking values that must stay text (matching the
# rest of this workbook's convention for this column group), so force
# a text number format before writing them, otherwise Excel-style
# auto-detection would coerce them into real numbers.
$q1.Range("D2:G2").NumberFormat = "@"
$q1.Range("D2").Value = "0.39"
$q1.Range("E2").Value = "80.03"
$q1.Range("F2").Value = "1.69"
$q1.Range("G2").Value = "0.0066"

$q1.Range("H2").Value = 5

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet at the end with the combined summary
# for both quarters.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)

$total.Range("A2").Value = 0
$q1.Range("A2").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$q1.Range("A2").Copy()
$total.Range("A3").PasteSpecial($xlPasteFormats)
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0
